$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '57.405.23'
$ws.Range("E2").Value = '  -0.44%  '

# Row 3
$ws.Range("D3").Value = '3.083.01'
$ws.Range("E3").Value = '  +0.59%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '514.84'
$ws.Range("E5").Value = '  -0.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.01'
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("E8").Value = '  -0.52%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.25'
$ws.Range("E9").Value = '  +0.26%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  -1.26%  '

# Row 11
$ws.Range("E11").Value = '  -1.44%  '

# Row 12
$ws.Range("D12").Value = '3.618.97'
$ws.Range("E12").Value = '  +0.95%  '

# Row 13
$ws.Range("E13").Value = '  +2.59%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.62'
$ws.Range("E14").Value = '  -4.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000163'
$ws.Range("E15").Value = '  -3.01%  '

# Row 16
$ws.Range("D16").Value = '57.534.06'
$ws.Range("E16").Value = '  -0.35%  '

# Row 17
$ws.Range("D17").Value = '3.089.40'
$ws.Range("E17").Value = '  +0.85%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.11'
$ws.Range("E18").Value = '  -1.59%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.08'
$ws.Range("E19").Value = '  -3.17%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  -0.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '335.02'
$ws.Range("E21").Value = '  +0.90%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.19%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.501'
$ws.Range("E23").Value = '  -1.47%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.68'
$ws.Range("E24").Value = '  +0.58%  '

# Row 25
$ws.Range("E25").Value = '  +3.91%  '

# Row 26
$ws.Range("E26").Value = '  +0.19%  '

# Row 27
$ws.Range("D27").Value = '0.0₃0922'
$ws.Range("E27").Value = '  +1.19%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.42'
$ws.Range("E28").Value = '  -5.14%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.10'
$ws.Range("E29").Value = '  -2.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.81'
$ws.Range("E30").Value = '  +0.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.87'
$ws.Range("E31").Value = '  -0.25%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.16'
$ws.Range("E32").Value = '  -4.56%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '153.93'
$ws.Range("E33").Value = '  -0.01%  '

# Row 34
$ws.Range("B34").Value = 'EnergySwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.58'
$ws.Range("E34").Value = '  +9.78%  '

# Row 35
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.53'
$ws.Range("E35").Value = '  -3.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.89'
$ws.Range("E36").Value = '  -0.80%  '

# Row 37
$ws.Range("E37").Value = '  -2.92%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0680'
$ws.Range("E38").Value = '  -0.28%  '

# Row 39
$ws.Range("D39").Value = '3.129.18'
$ws.Range("E39").Value = '  +1.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.78'
$ws.Range("E40").Value = '  -1.18%  '

# Row 41
$ws.Range("E41").Value = '  +0.33%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("E42").Value = '  -1.48%  '

# Row 43
$ws.Range("E43").Value = '  +0.16%  '

# Row 44
$ws.Range("D44").Value = '2.295.80'
$ws.Range("E44").Value = '  +3.88%  '

# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0253'
$ws.Range("E45").Value = '  +3.81%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.39'
$ws.Range("E46").Value = '  -1.14%  '

# Row 47
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.940'
$ws.Range("E47").Value = '  -1.37%  '

# Row 48
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.09'
$ws.Range("E48").Value = '  -0.91%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.88'
$ws.Range("E49").Value = '  -3.89%  '

# Row 50
$ws.Range("E50").Value = '  +1.18%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '246.95'
$ws.Range("E51").Value = '  +6.39%  '
